$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update email order in B3
$ws.Range("B3").Value = "comercial@elevadoreskorman.com.br;vendas@elevadoreskorman.com.br;korman@elevadoreskorman.com.br;"

# Row 4 becomes what used to be row 5 (coteibem)
$ws.Range("A4").Value = "https://coteibem.sindiconet.com.br/fornecedores/manutencao-elevadores/sp/sao-paulo"
$ws.Range("B4").Value = "contato@coteibem.com.br;"

# Delete old rows 5-8 (primac, framartel, continental, and old row5 duplicate shift)
$ws.Range("A5:B8").EntireRow.Delete()
